$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $val) {
    $r = $ws.Range($rangeAddr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "26.919.29"
$ws.Range("E2").Value = "  -0.42%  "
Set-TextValue "D3" "1.668.45"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  +0.03%  "
Set-TextValue "D6" "0.521"
$ws.Range("E6").Value = "  +2.14%  "
$ws.Range("E7").Value = "  +0.01%  "
Set-TextValue "D8" "0.0622"
$ws.Range("E8").Value = "  +1.22%  "
$ws.Range("E9").Value = "  +0.16%  "
Set-TextValue "D10" "20.32"
$ws.Range("E10").Value = "  +1.78%  "
$ws.Range("E11").Value = "  +2.62%  "
Set-TextValue "D12" "1.903.76"
$ws.Range("E12").Value = "  +0.91%  "
Set-TextValue "D13" "1.684.96"
$ws.Range("E13").Value = "  +1.81%  "
$ws.Range("E15").Value = "  +1.22%  "
Set-TextValue "D16" "65.65"
$ws.Range("E16").Value = "  +0.59%  "
Set-TextValue "D17" "26.905.87"
$ws.Range("E17").Value = "  -0.46%  "
Set-TextValue "D18" "235.57"
$ws.Range("E18").Value = "  -1.66%  "
$ws.Range("E19").Value = "  +1.66%  "
$ws.Range("E20").Value = "  +0.21%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("E22").Value = "  -0.53%  "
$ws.Range("E23").Value = "  -0.89%  "
$ws.Range("E24").Value = "  -3.30%  "
Set-TextValue "D25" "146.65"
$ws.Range("E25").Value = "  +0.41%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").Value = "  -1.34%  "
Set-TextValue "D28" "15.88"
$ws.Range("E28").Value = "  +0.33%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  -0.64%  "
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("E32").Value = "  +0.96%  "
Set-TextValue "D33" "1.447.20"
$ws.Range("E33").Value = "  -4.49%  "
Set-TextValue "D34" "3.14"
$ws.Range("E34").Value = "  +1.95%  "
$ws.Range("E35").Value = "  +2.74%  "
$ws.Range("E36").Value = "  -0.04%  "
Set-TextValue "D37" "0.585"
$ws.Range("E37").Value = "  +1.02%  "
$ws.Range("E38").Value = "  +1.56%  "
$ws.Range("E39").Value = "  +0.67%  "
Set-TextValue "D40" "5.73"
$ws.Range("E40").Value = "  -3.99%  "
$ws.Range("E42").Value = "  +8.82%  "
$ws.Range("E43").Value = "  +1.94%  "
Set-TextValue "D44" "66.07"
$ws.Range("E44").Value = "  +0.42%  "
Set-TextValue "D45" "1.809.72"
$ws.Range("E46").Value = "  +0.75%  "
Set-TextValue "D47" "90.72"
$ws.Range("E47").Value = "  +0.96%  "
$ws.Range("E48").Value = "  +1.00%  "
$ws.Range("E49").Value = "  +4.08%  "
Set-TextValue "D50" "0.0508"
Set-TextValue "D51" "7.59"
$ws.Range("E51").Value = "  -0.06%  "
